$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 duplicates row 2's data exactly (same match/venue/batsman stats
# row repeated in the source JSON export), so copy row 2 down to row 3.
# Using Copy/PasteSpecial (rather than re-typing the values) preserves
# the original "text" cell type for the numeric-looking values (runs,
# balls, 4s, 6s, strike rate), matching how row 2 itself is stored.
$ws.Range("A2:K2").Copy()
$ws.Range("A3:K3").PasteSpecial()
